# Update cached market-board price/profit figures on each job sheet
# (columns H-N) per the latest scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 359.53845
$ws.Range("I53").Value = 271.66666
$ws.Range("J53").Value = 434.85715
$ws.Range("K53").Value = 271.66666
$ws.Range("L53").Value = 434.85715
$ws.Range("M53").Value = 365.33334
$ws.Range("N53").Value = -1708.85715
$ws.Range("H103").Value = 1029467.9
$ws.Range("I103").Value = 609.5
$ws.Range("J103").Value = 3087184.8
$ws.Range("K103").Value = 1828.5
$ws.Range("L103").Value = 9261554.399999999
$ws.Range("M103").Value = -1242.5
$ws.Range("N103").Value = -9262726.399999999
$ws.Range("H132").Value = 226186.58
$ws.Range("I132").Value = 4044.2974
$ws.Range("J132").Value = 1253594.6
$ws.Range("K132").Value = 12132.8922
$ws.Range("L132").Value = 3760783.8
$ws.Range("M132").Value = -9602.8922
$ws.Range("N132").Value = -3765843.8
$ws.Range("H137").Value = 5813.8184
$ws.Range("I137").Value = 900.5
$ws.Range("J137").Value = 8621.429
$ws.Range("K137").Value = 2701.5
$ws.Range("L137").Value = 25864.287
$ws.Range("M137").Value = -151.5
$ws.Range("N137").Value = -30964.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3171.25
$ws.Range("I2").Value = 905.78125
$ws.Range("J2").Value = 9212.5
$ws.Range("K2").Value = 905.78125
$ws.Range("L2").Value = 9212.5
$ws.Range("M2").Value = -792.78125
$ws.Range("N2").Value = -9438.5
$ws.Range("H32").Value = 3355.86
$ws.Range("I32").Value = 2874.516
$ws.Range("J32").Value = 9750.857
$ws.Range("K32").Value = 2874.516
$ws.Range("L32").Value = 9750.857
$ws.Range("M32").Value = -2587.516
$ws.Range("N32").Value = -10324.857
$ws.Range("H63").Value = 1502.4166
$ws.Range("I63").Value = 1411.7273
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1411.7273
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -725.7273
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 1502.4166
$ws.Range("I66").Value = 1411.7273
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 7058.636500000001
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -3626.636500000001
$ws.Range("N66").Value = -19364
$ws.Range("H116").Value = 3171.25
$ws.Range("I116").Value = 905.78125
$ws.Range("J116").Value = 9212.5
$ws.Range("K116").Value = 905.78125
$ws.Range("L116").Value = 9212.5
$ws.Range("M116").Value = 1388.21875
$ws.Range("N116").Value = -13800.5
$ws.Range("H122").Value = 1149.7273
$ws.Range("I122").Value = 970
$ws.Range("J122").Value = 1464.25
$ws.Range("K122").Value = 2910
$ws.Range("L122").Value = 4392.75
$ws.Range("M122").Value = -460
$ws.Range("N122").Value = -9292.75
$ws.Range("H132").Value = 235481.73
$ws.Range("I132").Value = 51264.207
$ws.Range("J132").Value = 419699.25
$ws.Range("K132").Value = 153792.621
$ws.Range("L132").Value = 1259097.75
$ws.Range("M132").Value = -151262.621
$ws.Range("N132").Value = -1264157.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3171.25
$ws.Range("I3").Value = 905.78125
$ws.Range("J3").Value = 9212.5
$ws.Range("K3").Value = 905.78125
$ws.Range("L3").Value = 9212.5
$ws.Range("M3").Value = -791.78125
$ws.Range("N3").Value = -9440.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12057.714
$ws.Range("I31").Value = 8107.9287
$ws.Range("J31").Value = 19957.285
$ws.Range("K31").Value = 8107.9287
$ws.Range("L31").Value = 19957.285
$ws.Range("M31").Value = -7812.9287
$ws.Range("N31").Value = -20547.285
$ws.Range("H34").Value = 12057.714
$ws.Range("I34").Value = 8107.9287
$ws.Range("J34").Value = 19957.285
$ws.Range("K34").Value = 8107.9287
$ws.Range("L34").Value = 19957.285
$ws.Range("M34").Value = -7905.9287
$ws.Range("N34").Value = -20361.285
$ws.Range("H122").Value = 1232.0869
$ws.Range("I122").Value = 1308.2222
$ws.Range("J122").Value = 958
$ws.Range("K122").Value = 3924.6666
$ws.Range("L122").Value = 2874
$ws.Range("M122").Value = -1474.6666
$ws.Range("N122").Value = -7774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3019.3333
$ws.Range("I5").Value = 3019.3333
$ws.Range("K5").Value = 9057.999899999999
$ws.Range("M5").Value = -8945.999899999999
$ws.Range("H45").Value = 775
$ws.Range("I45").Value = 100
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 300
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = 232
$ws.Range("N45").Value = -4064
$ws.Range("H104").Value = 2792.6667
$ws.Range("J104").Value = 2792.6667
$ws.Range("L104").Value = 8378.000100000001
$ws.Range("N104").Value = -13620.0001
$ws.Range("H113").Value = 407.36365
$ws.Range("I113").Value = 370.57144
$ws.Range("J113").Value = 471.75
$ws.Range("K113").Value = 1111.71432
$ws.Range("L113").Value = 1415.25
$ws.Range("M113").Value = 1058.28568
$ws.Range("N113").Value = -5755.25
$ws.Range("H122").Value = 605.0909
$ws.Range("J122").Value = 1251.6
$ws.Range("L122").Value = 11264.4
$ws.Range("N122").Value = -16164.4
$ws.Range("H131").Value = 78947830
$ws.Range("I131").Value = 339.0909
$ws.Range("J131").Value = 187500620
$ws.Range("K131").Value = 1017.2727
$ws.Range("L131").Value = 562501860
$ws.Range("M131").Value = 4022.7273
$ws.Range("N131").Value = -562511940
$ws.Range("H135").Value = 3019.3333
$ws.Range("I135").Value = 3019.3333
$ws.Range("K135").Value = 27173.9997
$ws.Range("M135").Value = -24638.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 74504
$ws.Range("J29").Value = 74504
$ws.Range("L29").Value = 74504
$ws.Range("N29").Value = -75084
$ws.Range("H70").Value = 4148.857
$ws.Range("I70").Value = 3378.4
$ws.Range("J70").Value = 6075
$ws.Range("K70").Value = 3378.4
$ws.Range("L70").Value = 6075
$ws.Range("M70").Value = -3108.4
$ws.Range("N70").Value = -6615
$ws.Range("H73").Value = 4148.857
$ws.Range("I73").Value = 3378.4
$ws.Range("J73").Value = 6075
$ws.Range("K73").Value = 3378.4
$ws.Range("L73").Value = 6075
$ws.Range("M73").Value = -2442.4
$ws.Range("N73").Value = -7947
$ws.Range("H132").Value = 43857.582
$ws.Range("I132").Value = 1931.6666
$ws.Range("J132").Value = 169635.33
$ws.Range("K132").Value = 5794.9998
$ws.Range("L132").Value = 508905.99
$ws.Range("M132").Value = -3264.9998
$ws.Range("N132").Value = -513965.99

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2966.4707
$ws.Range("I61").Value = 2516.6667
$ws.Range("J61").Value = 3472.5
$ws.Range("K61").Value = 2516.6667
$ws.Range("L61").Value = 3472.5
$ws.Range("M61").Value = -2314.6667
$ws.Range("N61").Value = -3876.5
$ws.Range("H113").Value = 2966.4707
$ws.Range("I113").Value = 2516.6667
$ws.Range("J113").Value = 3472.5
$ws.Range("K113").Value = 2516.6667
$ws.Range("L113").Value = 3472.5
$ws.Range("M113").Value = -346.6667000000002
$ws.Range("N113").Value = -7812.5
$ws.Range("H132").Value = 155086.39
$ws.Range("I132").Value = 25445.861
$ws.Range("J132").Value = 347312
$ws.Range("K132").Value = 76337.583
$ws.Range("L132").Value = 1041936
$ws.Range("M132").Value = -73807.583
$ws.Range("N132").Value = -1046996
$ws.Range("H133").Value = 48142.6
$ws.Range("J133").Value = 48142.6
$ws.Range("L133").Value = 48142.6
$ws.Range("N133").Value = -53202.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 358.83334
$ws.Range("I113").Value = 373.27274
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 1119.81822
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1050.18178
$ws.Range("N113").Value = -4940
$ws.Range("H132").Value = 2268.127
$ws.Range("I132").Value = 429
$ws.Range("J132").Value = 5946.381
$ws.Range("K132").Value = 1287
$ws.Range("L132").Value = 17839.143
$ws.Range("M132").Value = 1243
$ws.Range("N132").Value = -22899.143
